# Update NASDAQ-100 ticker list: sector / sub-sector reclassifications
# as per GICS taxonomy refresh (2023-09-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - ADP
$ws.Range("C3").Value = "Industrials"
$ws.Range("D3").Value = "Human Resource & Employment Services"

# Row 4 - Airbnb
$ws.Range("D4").Value = "Hotels, Resorts & Cruise Lines"

# Row 8 - Amazon
$ws.Range("D8").Value = "Broadline Retail"

# Row 22 - Booking Holdings
$ws.Range("D22").Value = "Hotels, Resorts & Cruise Lines"

# Row 33 - Costco
$ws.Range("D33").Value = "Consumer Staples Merchandise Retail"

# Row 35 - CSX Corporation
$ws.Range("D35").Value = "Rail Transportation"

# Row 39 - Dollar Tree
$ws.Range("C39").Value = "Consumer Staples"
$ws.Range("D39").Value = "Consumer Staples Merchandise Retail"

# Row 40 - eBay
$ws.Range("D40").Value = "Broadline Retail"

# Row 55 - JD.com
$ws.Range("D55").Value = "Broadline Retail"

# Row 56 - Keurig Dr Pepper
$ws.Range("D56").Value = "Soft Drinks & Non-alcoholic Beverages"

# Row 64 - MercadoLibre
$ws.Range("D64").Value = "Broadline Retail"

# Row 71 - Monster Beverage
$ws.Range("D71").Value = "Soft Drinks & Non-alcoholic Beverages"

# Row 75 - O'Reilly Automotive
$ws.Range("D75").Value = "Automotive Retail"

# Row 76 - Old Dominion Freight Line
$ws.Range("D76").Value = "Cargo Ground Transportation"

# Row 78 - Paccar
$ws.Range("D78").Value = "Construction Machinery & Heavy Transportation Equipment"

# Row 80 - Paychex
$ws.Range("C80").Value = "Industrials"
$ws.Range("D80").Value = "Human Resource & Employment Services"

# Row 81 - PayPal
$ws.Range("C81").Value = "Financials"
$ws.Range("D81").Value = "Transaction & Payment Processing Services"

# Row 82 - PDD Holdings
$ws.Range("D82").Value = "Broadline Retail"

# Row 83 - PepsiCo
$ws.Range("D83").Value = "Soft Drinks & Non-alcoholic Beverages"
